$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.214.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.31%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.185.26"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -7.36%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.51"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.602"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.39%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.180.90"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -7.41%  "
$ws.Range("E10").Value = "  -5.05%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.64"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.83%  "
$ws.Range("E12").Value = "  -2.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.734.30"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -7.37%  "
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "27.67"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.91%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.179.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.43%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000163"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.184.92"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -7.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.69"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.74%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "353.15"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.62%  "
$ws.Range("E23").Value = "  +0.10%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.49"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.49%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.505"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.86%  "
$ws.Range("E26").Value = "  -1.68%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.26%  "
$ws.Range("E28").Value = "  -0.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.995"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.53%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.69"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.06%  "
$ws.Range("E31").Value = "  -0.06%  "
$ws.Range("E32").Value = "  -4.00%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.18"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.95%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.66"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.12%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "156.95"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.09%  "
$ws.Range("E37").Value = "  -5.30%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.809"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -7.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "25.96"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -8.86%  "
$ws.Range("E40").Value = "  -3.17%  "
$ws.Range("E41").Value = "  -3.85%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.670.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.16%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.17"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -5.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.03"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -6.18%  "
$ws.Range("B45").Value = "Hedera"
$ws.Range("C45").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0652"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.69%  "
$ws.Range("B46").Value = "Bittensor"
$ws.Range("C46").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "327.37"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "38.60"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.00%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.61%  "
$ws.Range("E49").Value = "  -5.81%  "
$ws.Range("E50").Value = "  -0.30%  "
$ws.Range("E51").Value = "  -0.03%  "
